$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V6")

# Insert 4 new rows at the top of the tracked-event list (before old row 3,
# "01 Start"). This pushes the old row 3 ("01 Start") down to row 7 and
# every subsequent row down by 4 as well.
$ws.Rows("3:6").Insert()

# The old "01 Start" row (now sitting at row 7) is being replaced by more
# granular tracking points, so remove it. Everything below shifts back up
# by one, giving a net shift of +3 rows overall.
$ws.Rows("7:7").Delete()

# Fill in the four new tracked events.
$ws.Range("A3").Value = "01 1st move"
$ws.Range("B3").Value = 529
$ws.Range("C3").Value = 736
$ws.Range("D3").Formula = "=IF(C3<>"""",IF(B3<>"""",C3-B3,""-""), ""-"")"

$ws.Range("A4").Value = "gold left 4"
$ws.Range("B4").Value = 600
$ws.Range("C4").Value = 808
$ws.Range("D4").Formula = "=IF(C4<>"""",IF(B4<>"""",C4-B4,""-""), ""-"")"

$ws.Range("A5").Value = "gold left 3"
$ws.Range("B5").Value = 628
$ws.Range("C5").Value = 840
$ws.Range("D5").Formula = "=IF(C5<>"""",IF(B5<>"""",C5-B5,""-""), ""-"")"

$ws.Range("A6").Value = "gold left 2"
$ws.Range("B6").Value = 663
$ws.Range("C6").Value = 872
$ws.Range("D6").Formula = "=IF(C6<>"""",IF(B6<>"""",C6-B6,""-""), ""-"")"

# Move the active selection to where the user was last working.
[void]$ws.Range("C7").Select()
